$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) "总计" sheet: insert a new data row for 2022-Q4 right under the header,
#    pushing the existing quarters down by one row.
# ---------------------------------------------------------------------------
$summary = $wb.Worksheets.Item("总计")

$summary.Rows.Item(2).Insert()

$summary.Cells.Item(2,1).Value = 0
$summary.Cells.Item(2,2).Value = "2022-Q4"
$summary.Cells.Item(2,3).Value = 3
$summary.Cells.Item(2,4).Value = 0.03

# Re-number the running index in column A for the rows that shifted down.
for ($r = 3; $r -le 8; $r++) {
    $summary.Cells.Item($r,1).Value = $r - 2
}

# The row-insert leaves stray formatting behind; restore the same look the
# other data rows have (bordered/bold index cell in column A, plain cells
# elsewhere).
$summary.Range("B2:D2").ClearFormats()
$summary.Cells.Item(3,1).Copy()
$summary.Cells.Item(2,1).PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# 2) Add a brand-new "2022-Q4" sheet (right after "总计", before "2022-Q3")
#    holding the fund-holdings detail for the new quarter.
# ---------------------------------------------------------------------------
# Duplicate an existing quarter sheet that already has the right shape (1
# header row + 3 data rows) so the new sheet inherits identical sheet-level
# properties (outline/page-setup), column layout and cell styles; its values
# get overwritten below.
$template = $wb.Worksheets.Item("2022-Q2")
$template.Copy($wb.Worksheets.Item(2))
$newSheet = $wb.Worksheets.Item(2)
$newSheet.Name = "2022-Q4"

# Columns B-G are textual in the source data (codes/names/percentages kept
# as strings, e.g. to preserve leading zeros), so force text format before
# assigning, then drop the format back to Normal to avoid a stray style.
$newSheet.Range("B2:G4").NumberFormat = "@"

$newSheet.Range("A2").Value = 0
$newSheet.Range("B2").Value = "004223"
$newSheet.Range("C2").Value = "金信多策略精选灵活配置混合"
$newSheet.Range("D2").Value = "0.31"
$newSheet.Range("E2").Value = "93.96"
$newSheet.Range("F2").Value = "4.83"
$newSheet.Range("G2").Value = "0.0150"
$newSheet.Range("H2").Value = 8

$newSheet.Range("A3").Value = 1
$newSheet.Range("B3").Value = "014246"
$newSheet.Range("C3").Value = "大摩现代服务业混合A"
$newSheet.Range("D3").Value = "0.17"
$newSheet.Range("E3").Value = "86.98"
$newSheet.Range("F3").Value = "6.17"
$newSheet.Range("G3").Value = "0.0105"
$newSheet.Range("H3").Value = 6

$newSheet.Range("A4").Value = 2
$newSheet.Range("B4").Value = "014247"
$newSheet.Range("C4").Value = "大摩现代服务业混合C"
$newSheet.Range("D4").Value = "0.06"
$newSheet.Range("E4").Value = "86.98"
$newSheet.Range("F4").Value = "6.17"
$newSheet.Range("G4").Value = "0.0037"
$newSheet.Range("H4").Value = 6

$newSheet.Range("B2:G4").Style = "Normal"

$summary.Activate()
